$d = $word.ActiveDocument

# The three styled runs ("some text" / "more text" / "even more text") live in
# the default footer (wdHeaderFooterPrimary) of the only section.
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)

$targets = "some text", "more text", "even more text"

foreach ($target in $targets) {
    $rng = $footer.Range.Duplicate
    $found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        # Re-assert the run's bold/italic/strike-through formatting so the
        # run properties are rewritten (POI 5.2.3 now emits on/off rather
        # than true/false for these toggle properties).
        $rng.Font.Bold = $true
        $rng.Font.Italic = $false
        $rng.Font.StrikeThrough = $false
    }
}
